$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Make the header row (row 1) bold ---
$ws.Range("A1:F1").Font.Bold = $true

# --- Widen column C to fit the new, longer "Module" values ---
$ws.Columns.Item(3).ColumnWidth = 26.14

# --- New future-enhancement rows ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Multiple Payment for one invoice"
$ws.Range("C6").Value = "Invoice"
$ws.Range("D6").Value = "Open"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Revert check out in case of operator mistake"
$ws.Range("C7").Value = "Checkin"
$ws.Range("D7").Value = "Open"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Recycle bin"
$ws.Range("C8").Value = "Navigator and all component"
$ws.Range("D8").Value = "Open"

# --- Match the new selection/active cell left by the edit ---
$ws.Range("B8").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
